# Weekly update: a new Ciboulette price record (week of 2023-10-20) was
# added as the new first data row (row 13) of the sheet. Every existing
# record from the old row 13 down to the old row 32 shifts down by one
# row (old row 13 -> row 14, ..., old row 32 -> row 33), and the new
# record's values are written into the freed-up row 13.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at row 13, pushing rows 13:32 down to 14:33.
$ws.Rows.Item(13).Insert()

# Populate the new row 13 with the new weekly record.
$ws.Range("A13").Value = 7
$ws.Range("B13").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C13").Value = "Ñuble"
$ws.Range("D13").Value = 45219
$ws.Range("E13").Value = 16
$ws.Range("F13").Value = 100112039
$ws.Range("G13").Value = "Ciboulette"
$ws.Range("H13").Value = "Sin especificar"
$ws.Range("I13").Value = "Primera"
$ws.Range("J13").Value = 300
$ws.Range("K13").Value = 2000
$ws.Range("L13").Value = 2500
$ws.Range("M13").Value = 2250
$ws.Range("N13").Value = "$/docena de atados"
$ws.Range("O13").Value = "Región Metropolitana"
$ws.Range("P13").Value = 750
$ws.Range("Q13").Value = 3
$ws.Range("R13").Value = "Hortaliza"
